# Mise à jour de l'application
# Adds 9 new training-load entries (row 561-569) for session dated 2025-11-13
# (Excel serial 45974), plus a brand new "Cheville (coup)" injury-location value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Style templates already present on the sheet:
#  - row 559 -> "Localisation douleur" (col G) filled with text (style s="1")
#  - row 560 -> "Localisation douleur" (col G) left empty          (style s="2")
$styleWithLocation   = $ws.Range("A559:I559")
$styleWithoutLocation = $ws.Range("A560:I560")

$dateSerial = 45974   # 2025-11-13

# row, player, volume, intensite, fatigue, douleur, localisation, plaisir
$newRows = @(
    @(561, "Yoan Zouma",       70, 4, 4, 6, "Ischio",           4),
    @(562, "Maé Clavel",       70, 6, 6, 0, $null,              6),
    @(563, "Karim Belmahi",    70, 6, 8, 0, $null,              10),
    @(564, "Naim Ighbane",     70, 6, 7, 8, "Genou",             3),
    @(565, "Omar Benyounes",   70, 5, 6, 2, "Cheville (coup)",   7),
    @(566, "Ilyes Boughanmi",  70, 6, 6, 2, "Genou",             8),
    @(567, "Levy Ndoutoume",   70, 7, 7, 1, "Ischio",            7),
    @(568, "Amir Etien",       70, 5, 5, 3, "Ischio",            2),
    @(569, "Sofiane Belle",    70, 5, 7, 2, "Dos",               7)
)

foreach ($row in $newRows) {
    $r          = $row[0]
    $player     = $row[1]
    $volume     = $row[2]
    $intensite  = $row[3]
    $fatigue    = $row[4]
    $douleur    = $row[5]
    $location   = $row[6]
    $plaisir    = $row[7]

    $destRow = $ws.Range("A" + $r + ":I" + $r)

    if ($location) {
        $styleWithLocation.Copy($destRow) | Out-Null
    } else {
        $styleWithoutLocation.Copy($destRow) | Out-Null
    }

    $ws.Range("A$r").Value = $dateSerial
    $ws.Range("B$r").Value = $player
    $ws.Range("C$r").Value = $volume
    $ws.Range("D$r").Value = $intensite
    $ws.Range("E$r").Value = $fatigue
    $ws.Range("F$r").Value = $douleur

    if ($location) {
        $ws.Range("G$r").Value = $location
    } else {
        $ws.Range("G$r").ClearContents() | Out-Null
    }

    $ws.Range("H$r").Value = $plaisir
    $ws.Range("I$r").Formula = "=C$r*D$r"
}

# Update the visible selection to match the latest entry.
$ws.Range("K564").Select() | Out-Null
